$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Audio" / "audio.mp3" column (was column A) - the rest of the
# data (columns B:G) shifts left to become the new A:F.
$ws.Columns("A").Delete()

# Reflect the post-delete selection state recorded in the saved file
# (Excel leaves the selection spanning the column that used to be selected).
$ws.Range("A1:A1048576").Select() | Out-Null
